$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set D8 to "V" (same marker used by the rest of column D) - row 8 was
# previously missing this cell entirely.
$ws.Range("D8").Value = "V"

# Move the active selection on the sheet from G17 to F10.
$ws.Range("F10").Select()
